# Regenerate orders with updated distance/sizes.
# Applies token-level text substitutions across every string cell in the
# used range of the active sheet:
#   D64 -> D69
#   D80 -> D86
#   D51 -> D55
#   S30 -> S31
# This mirrors the shared-strings table change in the target workbook,
# where Condition/Filename_Left/Filename_Right/Distance/Size columns all
# encode the distance/size tokens inside their text values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count()
$colCount = $used.Columns.Count()

$changed = 0

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()

        if ($val -is [string]) {
            $newVal = $val.Replace("D64", "D69").Replace("D80", "D86").Replace("D51", "D55").Replace("S30", "S31")

            if ($newVal -ne $val) {
                $cell.Value = $newVal
                $changed = $changed + 1
            }
        }
    }
}

Write-Host ("cells changed=" + $changed)
